# Final Project. Neighborhoods are totes working!!!
#
# Splits the sentence about neighborhood selection/highlighting into three
# runs (same formatting) while rewording it from:
#   "...that neighborhood is highlighted on the map on zoomed in."
# to:
#   "...that neighborhood is highlighted on the map. I also zoom in on the
#    neighborhood that is selected from the dropdown."

$d = $word.ActiveDocument

$oldText = "When the user selects a particular neighborhood, that neighborhood is highlighted on the map on zoomed in."

$part1 = "When the user selects a particular neighborhood, that neighbo"
$part2 = "rhood is highlighted on the map. I a"
$part3 = "lso zoom in on the neighborhood that is selected from the dropdown."

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence to edit."
}

$start = $rng.Start
$end = $rng.End

# Replace the whole sentence with the new combined text first (keeps the
# original run / formatting for the whole span).
$rAll = $d.Range($start, $end)
$rAll.Text = $part1 + $part2 + $part3

# Now force the 2nd and 3rd segments into their own runs (identical
# formatting to the first) by briefly toggling a character property on each
# sub-range; this creates a genuine run boundary without changing the
# resulting visible formatting.
$p2Start = $start + $part1.Length
$p2End = $p2Start + $part2.Length
$r2 = $d.Range($p2Start, $p2End)
$r2.Bold = 1
$r2.Bold = 0

$p3Start = $p2End
$p3End = $p3Start + $part3.Length
$r3 = $d.Range($p3Start, $p3End)
$r3.Bold = 1
$r3.Bold = 0

Write-Output "Replaced sentence; final text: $($d.Range($start, $p3End).Text)"
